# Auto-generated edit script applying numeric "wish count" (F column) refreshes
# plus the one full event replacement (row 16 on "全部类型") described in the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1348
$ws.Range("F4").Value = 2104
$ws.Range("F5").Value = 796
$ws.Range("F6").Value = 1532
$ws.Range("F7").Value = 97470
$ws.Range("F8").Value = 40547
$ws.Range("F10").Value = 104
$ws.Range("F11").Value = 625
$ws.Range("G11").Value = '不可售'
$ws.Range("F12").Value = 131
$ws.Range("F14").Value = 1426
$ws.Range("F18").Value = 868
$ws.Range("F19").Value = 5513
$ws.Range("F20").Value = 372
$ws.Range("F21").Value = 1026
$ws.Range("F22").Value = 2613
$ws.Range("F23").Value = 6164
$ws.Range("F24").Value = 152
$ws.Range("F25").Value = 1120
$ws.Range("F26").Value = 657
$ws.Range("F27").Value = 78
$ws.Range("F29").Value = 1081
$ws.Range("F31").Value = 50
$ws.Range("F35").Value = 816
$ws.Range("F37").Value = 72
$ws.Range("F42").Value = 139
$ws.Range("F43").Value = 37
$ws.Range("F44").Value = 163
$ws.Range("F45").Value = 1122
$ws.Range("F46").Value = 26
$ws.Range("F47").Value = 53

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 792
$ws.Range("F10").Value = 39
$ws.Range("F22").Value = 164
$ws.Range("F23").Value = 120
$ws.Range("F27").Value = 44
$ws.Range("F29").Value = 105
$ws.Range("F30").Value = 321
$ws.Range("F31").Value = 922
$ws.Range("F32").Value = 530
$ws.Range("F34").Value = 59
$ws.Range("F37").Value = 87
$ws.Range("F40").Value = 34

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 688
$ws.Range("F5").Value = 807
$ws.Range("F6").Value = 496

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 688
$ws.Range("F4").Value = 1348
$ws.Range("F5").Value = 807
$ws.Range("F8").Value = 496
$ws.Range("F13").Value = 1532
$ws.Range("F15").Value = 104
$ws.Range("C16").Value = '上海·多厨狂喜动漫展1.0'
$ws.Range("D16").Value = '澳门路168号 月星家居（澳门路）'
$ws.Range("E16").Value = '2024.05.03 10:00-05.04 16:00'
$ws.Range("F16").Value = 131
$ws.Range("G16").Value = 59
$ws.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=83932'
$ws.Range("I16").Value = '//i0.hdslb.com/bfs/openplatform/202404/721HW21G1712565123950.jpeg'
$ws.Range("F18").Value = 1426
$ws.Range("F21").Value = 39
$ws.Range("F22").Value = 5513
$ws.Range("F23").Value = 372
$ws.Range("F24").Value = 1026
$ws.Range("F25").Value = 2613
$ws.Range("F27").Value = 6164
$ws.Range("F28").Value = 152
$ws.Range("F29").Value = 1120
$ws.Range("F31").Value = 657
$ws.Range("F32").Value = 78
$ws.Range("F33").Value = 1081
$ws.Range("F35").Value = 816
$ws.Range("F37").Value = 72
$ws.Range("F41").Value = 922
$ws.Range("F42").Value = 530
$ws.Range("F43").Value = 139
$ws.Range("F44").Value = 59
$ws.Range("F45").Value = 163
$ws.Range("F46").Value = 87
$ws.Range("F49").Value = 53
$ws.Range("F50").Value = 34
